$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New trading log rows to append (rows 4-9), mirroring the existing
# timestamp/action/token/signal_type/price/.../status/error_message layout.
$rows = @(
    @{ A = "2025-09-19T19:26:56.251149"; B = "TRADING_ATTEMPT"; C = "NEAR"; D = "UNKNOWN"; E = 3.115658833304698;  K = "ATTEMPT"; L = "Attempting trade 1/3" },
    @{ A = "2025-09-19T19:26:57.769743"; B = "POSITION_FAILED"; C = "NEAR"; D = "UNKNOWN"; E = $null;             K = "FAILED";  L = "Trade execution failed for trade 1" },
    @{ A = "2025-09-19T19:26:57.782292"; B = "TRADING_ATTEMPT"; C = "SUI";  D = "UNKNOWN"; E = 3.66429971981016;   K = "ATTEMPT"; L = "Attempting trade 2/3" },
    @{ A = "2025-09-19T19:26:59.211472"; B = "POSITION_FAILED"; C = "SUI";  D = "UNKNOWN"; E = $null;             K = "FAILED";  L = "Trade execution failed for trade 2" },
    @{ A = "2025-09-19T19:26:59.224981"; B = "TRADING_ATTEMPT"; C = "ADA";  D = "UNKNOWN"; E = 0.8963706973452196; K = "ATTEMPT"; L = "Attempting trade 3/3" },
    @{ A = "2025-09-19T19:27:00.824559"; B = "POSITION_FAILED"; C = "ADA";  D = "UNKNOWN"; E = $null;             K = "FAILED";  L = "Trade execution failed for trade 3" }
)

$rowIndex = 4
foreach ($r in $rows) {
    $ws.Cells.Item($rowIndex, 1).Value = $r.A
    $ws.Cells.Item($rowIndex, 2).Value = $r.B
    $ws.Cells.Item($rowIndex, 3).Value = $r.C
    $ws.Cells.Item($rowIndex, 4).Value = $r.D
    if ($null -ne $r.E) {
        $ws.Cells.Item($rowIndex, 5).Value = $r.E
    }
    $ws.Cells.Item($rowIndex, 11).Value = $r.K
    $ws.Cells.Item($rowIndex, 12).Value = $r.L
    $rowIndex++
}
